$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 283.74
$ws.Range("I15").Value = 283.74
$ws.Range("K15").Value = 851.22
$ws.Range("M15").Value = -682.22
$ws.Range("H100").Value = 13889778
$ws.Range("I100").Value = 13889778
$ws.Range("K100").Value = 13889778
$ws.Range("M100").Value = -13889237
$ws.Range("H137").Value = 1770.7106
$ws.Range("I137").Value = 1322.5143
$ws.Range("J137").Value = 6999.6665
$ws.Range("K137").Value = 3967.5429
$ws.Range("L137").Value = 20998.9995
$ws.Range("M137").Value = -1417.5429
$ws.Range("N137").Value = -26098.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 101150
$ws.Range("I2").Value = 134441.14
$ws.Range("J2").Value = 1276.6
$ws.Range("K2").Value = 134441.14
$ws.Range("L2").Value = 1276.6
$ws.Range("M2").Value = -134328.14
$ws.Range("N2").Value = -1502.6
$ws.Range("H32").Value = 12200.809
$ws.Range("I32").Value = 8729.102000000001
$ws.Range("J32").Value = 21782.72
$ws.Range("K32").Value = 8729.102000000001
$ws.Range("L32").Value = 21782.72
$ws.Range("M32").Value = -8442.102000000001
$ws.Range("N32").Value = -22356.72
$ws.Range("H61").Value = 170842.56
$ws.Range("I61").Value = 4806
$ws.Range("J61").Value = 479196.2
$ws.Range("K61").Value = 4806
$ws.Range("L61").Value = 479196.2
$ws.Range("M61").Value = -4594
$ws.Range("N61").Value = -479620.2
$ws.Range("H74").Value = 7354176
$ws.Range("I74").Value = 994.6585
$ws.Range("K74").Value = 994.6585
$ws.Range("M74").Value = -120.6585
$ws.Range("H77").Value = 7354176
$ws.Range("I77").Value = 994.6585
$ws.Range("K77").Value = 4973.2925
$ws.Range("M77").Value = -605.2924999999996
$ws.Range("H116").Value = 101150
$ws.Range("I116").Value = 134441.14
$ws.Range("J116").Value = 1276.6
$ws.Range("K116").Value = 134441.14
$ws.Range("L116").Value = 1276.6
$ws.Range("M116").Value = -132147.14
$ws.Range("N116").Value = -5864.6
$ws.Range("H132").Value = 1757375.4
$ws.Range("I132").Value = 2147.4211
$ws.Range("K132").Value = 6442.263300000001
$ws.Range("M132").Value = -3912.263300000001
$ws.Range("H136").Value = 170842.56
$ws.Range("I136").Value = 4806
$ws.Range("J136").Value = 479196.2
$ws.Range("K136").Value = 14418
$ws.Range("L136").Value = 1437588.6
$ws.Range("M136").Value = -11868
$ws.Range("N136").Value = -1442688.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 101150
$ws.Range("I3").Value = 134441.14
$ws.Range("J3").Value = 1276.6
$ws.Range("K3").Value = 134441.14
$ws.Range("L3").Value = 1276.6
$ws.Range("M3").Value = -134327.14
$ws.Range("N3").Value = -1504.6
$ws.Range("H94").Value = 1357.3226
$ws.Range("I94").Value = 940.3182
$ws.Range("J94").Value = 2376.6667
$ws.Range("K94").Value = 940.3182
$ws.Range("L94").Value = 2376.6667
$ws.Range("M94").Value = -489.3182
$ws.Range("N94").Value = -3278.6667
$ws.Range("H122").Value = 67761.664
$ws.Range("J122").Value = 67761.664
$ws.Range("L122").Value = 67761.664
$ws.Range("N122").Value = -77561.664
$ws.Range("H134").Value = 32915.164
$ws.Range("I134").Value = 6983.7036
$ws.Range("J134").Value = 102930.1
$ws.Range("K134").Value = 20951.1108
$ws.Range("L134").Value = 308790.3
$ws.Range("M134").Value = -18416.1108
$ws.Range("N134").Value = -313860.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2067712.8
$ws.Range("I5").Value = 450.57144
$ws.Range("J5").Value = 3955213
$ws.Range("K5").Value = 1351.71432
$ws.Range("L5").Value = 11865639
$ws.Range("M5").Value = -1239.71432
$ws.Range("N5").Value = -11865863
$ws.Range("H131").Value = 3227104.2
$ws.Range("I131").Value = 8333774
$ws.Range("J131").Value = 1839.2632
$ws.Range("K131").Value = 25001322
$ws.Range("L131").Value = 5517.7896
$ws.Range("M131").Value = -24996282
$ws.Range("N131").Value = -15597.7896
$ws.Range("H135").Value = 2067712.8
$ws.Range("I135").Value = 450.57144
$ws.Range("J135").Value = 3955213
$ws.Range("K135").Value = 4055.14296
$ws.Range("L135").Value = 35596917
$ws.Range("M135").Value = -1520.14296
$ws.Range("N135").Value = -35601987

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8444.0625
$ws.Range("J80").Value = 2444.4443
$ws.Range("L80").Value = 2444.4443
$ws.Range("N80").Value = -4440.4443
$ws.Range("H83").Value = 8444.0625
$ws.Range("J83").Value = 2444.4443
$ws.Range("L83").Value = 12222.2215
$ws.Range("N83").Value = -22206.2215
$ws.Range("H113").Value = 45455590
$ws.Range("I113").Value = 83334130
$ws.Range("K113").Value = 83334130
$ws.Range("M113").Value = -83331960

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2967.2666
$ws.Range("I7").Value = 1977.1111
$ws.Range("J7").Value = 4452.5
$ws.Range("K7").Value = 1977.1111
$ws.Range("L7").Value = 4452.5
$ws.Range("M7").Value = -1865.1111
$ws.Range("N7").Value = -4676.5
$ws.Range("H55").Value = 133.55556
$ws.Range("I55").Value = 155
$ws.Range("J55").Value = 90.666664
$ws.Range("K55").Value = 155
$ws.Range("L55").Value = 90.666664
$ws.Range("M55").Value = 18
$ws.Range("N55").Value = -436.666664
$ws.Range("H61").Value = 2917.5
$ws.Range("I61").Value = 2567.8572
$ws.Range("J61").Value = 3733.3333
$ws.Range("K61").Value = 2567.8572
$ws.Range("L61").Value = 3733.3333
$ws.Range("M61").Value = -2365.8572
$ws.Range("N61").Value = -4137.3333
$ws.Range("H113").Value = 2917.5
$ws.Range("I113").Value = 2567.8572
$ws.Range("J113").Value = 3733.3333
$ws.Range("K113").Value = 2567.8572
$ws.Range("L113").Value = 3733.3333
$ws.Range("M113").Value = -397.8571999999999
$ws.Range("N113").Value = -8073.3333
$ws.Range("H126").Value = 2967.2666
$ws.Range("I126").Value = 1977.1111
$ws.Range("J126").Value = 4452.5
$ws.Range("K126").Value = 5931.3333
$ws.Range("L126").Value = 13357.5
$ws.Range("M126").Value = -3461.3333
$ws.Range("N126").Value = -18297.5
$ws.Range("H132").Value = 17248496
$ws.Range("I132").Value = 27787356
$ws.Range("J132").Value = 3086.7273
$ws.Range("K132").Value = 83362068
$ws.Range("L132").Value = 9260.1819
$ws.Range("M132").Value = -83359538
$ws.Range("N132").Value = -14320.1819
$ws.Range("H136").Value = 14265.869
$ws.Range("I136").Value = 14227.8
$ws.Range("J136").Value = 14295.154
$ws.Range("K136").Value = 42683.39999999999
$ws.Range("L136").Value = 42885.462
$ws.Range("M136").Value = -40133.39999999999
$ws.Range("N136").Value = -47985.462
$ws.Range("H139").Value = 69778.75
$ws.Range("J139").Value = 69778.75
$ws.Range("L139").Value = 69778.75
$ws.Range("N139").Value = -80058.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 29299
$ws.Range("J70").Value = 29333
$ws.Range("L70").Value = 29333
$ws.Range("N70").Value = -29963
$ws.Range("H73").Value = 29299
$ws.Range("J73").Value = 29333
$ws.Range("L73").Value = 29333
$ws.Range("N73").Value = -31517
$ws.Range("H113").Value = 1183
$ws.Range("I113").Value = 1268.375
$ws.Range("J113").Value = 500
$ws.Range("K113").Value = 3805.125
$ws.Range("L113").Value = 500
$ws.Range("M113").Value = -1635.125
$ws.Range("N113").Value = -5840
